$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple Taxonsorteringsordning (column B) renumbering ---
# 57880 -> 57884 on rows 2,3,4,8,9,11,12
$ws.Range("B2").Value = 57884
$ws.Range("B3").Value = 57884
$ws.Range("B4").Value = 57884
$ws.Range("B8").Value = 57884
$ws.Range("B9").Value = 57884
$ws.Range("B11").Value = 57884
$ws.Range("B12").Value = 57884

# 79239 -> 79243 on rows 5,10
$ws.Range("B5").Value = 79243
$ws.Range("B10").Value = 79243

# --- Rows 6 and 7 are swapped in full (every column except A/B keeps the
# ---  other row's content; A swaps directly, B swaps and is bumped by +4
# ---  like the rest of the column). Capture the "before" values first so
# ---  the write-back doesn't clobber data we still need to read.

$oldA6 = $ws.Range("A6").Value2
$oldB6 = $ws.Range("B6").Value2
$oldD6 = $ws.Range("D6").Value2
$oldE6 = $ws.Range("E6").Value2
$oldF6 = $ws.Range("F6").Value2
$oldG6 = $ws.Range("G6").Value2
$oldH6 = $ws.Range("H6").Value2
$oldP6 = $ws.Range("P6").Value2
$oldQ6 = $ws.Range("Q6").Value2
$oldR6 = $ws.Range("R6").Value2
$oldZ6 = $ws.Range("Z6").Value2
$oldAB6 = $ws.Range("AB6").Value2
$oldAC6 = $ws.Range("AC6").Value2

$oldA7 = $ws.Range("A7").Value2
$oldB7 = $ws.Range("B7").Value2
$oldD7 = $ws.Range("D7").Value2
$oldE7 = $ws.Range("E7").Value2
$oldF7 = $ws.Range("F7").Value2
$oldG7 = $ws.Range("G7").Value2
$oldH7 = $ws.Range("H7").Value2
$oldP7 = $ws.Range("P7").Value2
$oldQ7 = $ws.Range("Q7").Value2
$oldR7 = $ws.Range("R7").Value2
$oldZ7 = $ws.Range("Z7").Value2
$oldAB7 = $ws.Range("AB7").Value2
$oldAC7 = $ws.Range("AC7").Value2

# Row 6 becomes what row 7 used to hold
$ws.Range("A6").Value = $oldA7
$ws.Range("B6").Value = ($oldB7 + 4)
$ws.Range("D6").Value = $oldD7
$ws.Range("E6").Value = $oldE7
$ws.Range("F6").Value = $oldF7
$ws.Range("G6").Value = $oldG7
$ws.Range("H6").Value = $oldH7
$ws.Range("P6").Value = $oldP7
$ws.Range("Q6").Value = $oldQ7
$ws.Range("R6").Value = $oldR7
$ws.Range("Z6").Value = $oldZ7
$ws.Range("AB6").Value = $oldAB7
$ws.Range("AC6").Value = $oldAC7

# Row 7 becomes what row 6 used to hold
$ws.Range("A7").Value = $oldA6
$ws.Range("B7").Value = ($oldB6 + 4)
$ws.Range("D7").Value = $oldD6
$ws.Range("E7").Value = $oldE6
$ws.Range("F7").Value = $oldF6
$ws.Range("G7").Value = $oldG6
$ws.Range("H7").Value = $oldH6
$ws.Range("P7").Value = $oldP6
$ws.Range("Q7").Value = $oldQ6
$ws.Range("R7").Value = $oldR6
$ws.Range("Z7").Value = $oldZ6
$ws.Range("AB7").Value = $oldAB6
$ws.Range("AC7").Value = $oldAC6
